$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "duration"
$ws.Range("C12").Value = "4 [1,15]"
$ws.Range("D12").Value = "0 (0%)"
$ws.Range("E12").Value = "10 [5,20]"
$ws.Range("F12").Value = "0 (0%)"
